$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep the Price column as text (values like "45.779.25" / "98.66" must not
# be auto-converted to numbers by Excel's COM layer).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '45.779.25'
$ws.Range("E2").Value = '  -1.10%  '

$ws.Range("D3").Value = '2.596.65'
$ws.Range("E3").Value = '  -0.69%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").Value = '308.40'
$ws.Range("E5").Value = '  -1.55%  '

$ws.Range("D6").Value = '98.66'
$ws.Range("E6").Value = '  -1.91%  '

$ws.Range("D7").Value = '0.594'
$ws.Range("E7").Value = '  -0.76%  '

$ws.Range("D8").Value = '1.00'

$ws.Range("D9").Value = '0.575'
$ws.Range("E9").Value = '  -1.47%  '

$ws.Range("D10").Value = '38.51'
$ws.Range("E10").Value = '  -0.64%  '

$ws.Range("B11").Value = 'OKB'
$ws.Range("C11").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D11").Value = '53.99'
$ws.Range("E11").Value = '  -0.65%  '

$ws.Range("B12").Value = 'Dogecoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D12").Value = '0.0838'
$ws.Range("E12").Value = '  -0.57%  '

$ws.Range("D13").Value = '8.01'
$ws.Range("E13").Value = '  -4.07%  '

$ws.Range("D14").Value = '2.997.15'
$ws.Range("E14").Value = '  -0.56%  '

$ws.Range("E15").Value = '  +0.68%  '

$ws.Range("D16").Value = '2.603.06'
$ws.Range("E16").Value = '  -0.64%  '

$ws.Range("D17").Value = '0.908'
$ws.Range("E17").Value = '  -0.34%  '

$ws.Range("D18").Value = '14.71'
$ws.Range("E18").Value = '  -2.00%  '

$ws.Range("D19").Value = '45.871.19'
$ws.Range("E19").Value = '  -1.39%  '

$ws.Range("E20").Value = '  -1.17%  '

$ws.Range("D21").Value = '6.67'
$ws.Range("E21").Value = '  -1.26%  '

$ws.Range("D22").Value = '12.54'
$ws.Range("E22").Value = '  -5.34%  '

$ws.Range("D23").Value = '286.06'
$ws.Range("E23").Value = '  +11.83%  '

$ws.Range("D24").Value = '73.04'
$ws.Range("E24").Value = '  +2.50%  '

$ws.Range("E25").Value = '  -2.47%  '

$ws.Range("E26").Value = '  +0.43%  '

$ws.Range("D27").Value = '29.04'
$ws.Range("E27").Value = '  +2.70%  '

$ws.Range("E28").Value = '  -0.02%  '

$ws.Range("D29").Value = '4.06'
$ws.Range("E29").Value = '  +0.65%  '

$ws.Range("D30").Value = '10.61'
$ws.Range("E30").Value = '  -0.12%  '

$ws.Range("D31").Value = '38.25'
$ws.Range("E31").Value = '  -4.81%  '

$ws.Range("E32").Value = '  -3.10%  '

$ws.Range("D33").Value = '6.25'
$ws.Range("E33").Value = '  +1.14%  '

$ws.Range("D34").Value = '3.60'
$ws.Range("E34").Value = '  -3.39%  '

$ws.Range("D35").Value = '157.77'
$ws.Range("E35").Value = '  +2.94%  '

$ws.Range("E36").Value = '  -3.04%  '

$ws.Range("D37").Value = '2.81'
$ws.Range("E37").Value = '  -2.90%  '

$ws.Range("D38").Value = '0.0829'
$ws.Range("E38").Value = '  -0.96%  '

$ws.Range("D39").Value = '0.122'
$ws.Range("E39").Value = '  +2.73%  '

$ws.Range("D40").Value = '0.122'
$ws.Range("E40").Value = '  +0.14%  '

$ws.Range("D41").Value = '15.64'
$ws.Range("E41").Value = '  -8.87%  '

$ws.Range("E42").Value = '  -0.10%  '

$ws.Range("E43").Value = '  -2.66%  '

$ws.Range("D44").Value = '3.98'
$ws.Range("E44").Value = '  -5.63%  '

$ws.Range("D45").Value = '21.24'
$ws.Range("E45").Value = '  +0.68%  '

$ws.Range("D46").Value = '2.102.89'
$ws.Range("E46").Value = '  +2.76%  '

$ws.Range("D47").Value = '0.999'
$ws.Range("E47").Value = '  +0.03%  '

$ws.Range("D48").Value = '94.16'
$ws.Range("E48").Value = '  +2.93%  '

$ws.Range("D49").Value = '9.22'
$ws.Range("E49").Value = '  -0.79%  '

$ws.Range("D50").Value = '108.28'
$ws.Range("E50").Value = '  -1.96%  '

$ws.Range("D51").Value = '2.852.55'
$ws.Range("E51").Value = '  -0.85%  '
